$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[53.14777310142459, 74.5921791779116]"
$ws.Range("N2").Value = [double]"1.332267629550188e-15"
$ws.Range("O2").Value = [double]"1.332267629550188e-15"
$ws.Range("U2").Value = "[43.098336248067596, 56.90188764318687]"
$ws.Range("M3").Value = "[53.07741293239641, 75.49464457501804]"
$ws.Range("N3").Value = [double]"4.662936703425657e-15"
$ws.Range("O3").Value = [double]"4.662936703425657e-15"
$ws.Range("U3").Value = "[41.96183359681138, 55.209697740524525]"
$ws.Range("M4").Value = "[52.399766137121645, 76.13739701406956]"
$ws.Range("N4").Value = [double]"3.197442310920451e-14"
$ws.Range("O4").Value = [double]"3.197442310920451e-14"
$ws.Range("Q4").Value = "[1.8428161111147343, 2.220184598032427]"
$ws.Range("U4").Value = "[44.51835659427349, 57.58220901451525]"
$ws.Range("Y4").Value = [double]"14.99573573573585"
$ws.Range("Z4").Value = [double]"16.38852852852866"
$ws.Range("M5").Value = "[51.582113342013784, 76.54868206491851]"
$ws.Range("N5").Value = [double]"1.825206652483757e-13"
$ws.Range("O5").Value = [double]"1.825206652483757e-13"
$ws.Range("U5").Value = "[42.87515256901382, 55.73160650913795]"
$ws.Range("M6").Value = "[51.08086380268998, 76.05601861198937]"
$ws.Range("N6").Value = [double]"2.364775042451583e-13"
$ws.Range("O6").Value = [double]"2.364775042451583e-13"
$ws.Range("U6").Value = "[43.266828036950486, 56.07657796208319]"
$ws.Range("M7").Value = "[51.64426172206817, 74.9737189671556]"
$ws.Range("N7").Value = [double]"2.97539770599542e-14"
$ws.Range("O7").Value = [double]"2.97539770599542e-14"
$ws.Range("U7").Value = "[43.39725752129752, 56.22740881366073]"
$ws.Range("M8").Value = "[51.95048350600315, 73.8849956797687]"
$ws.Range("N8").Value = [double]"4.662936703425657e-15"
$ws.Range("O8").Value = [double]"4.662936703425657e-15"
$ws.Range("Q8").Value = "[2.7107636310254275, 3.08813211794312]"
$ws.Range("U8").Value = "[43.55040145407687, 56.43531976332853]"
$ws.Range("Y8").Value = [double]"11.79231231231241"
$ws.Range("Z8").Value = [double]"13.18510510510521"
$ws.Range("M9").Value = "[52.36192544180358, 71.91895337287757]"
$ws.Range("N9").Value = [double]"0"
$ws.Range("O9").Value = [double]"0"
$ws.Range("U9").Value = "[42.9956065947495, 55.830204328551574]"
$ws.Range("M10").Value = "[52.745301513467325, 72.92474356988734]"
$ws.Range("N10").Value = [double]"2.220446049250313e-16"
$ws.Range("O10").Value = [double]"2.220446049250313e-16"
$ws.Range("U10").Value = "[44.9458167579161, 57.68321276229904]"
$ws.Range("M11").Value = "[52.15533663150874, 75.91184040506133]"
$ws.Range("N11").Value = [double]"3.708144902248023e-14"
$ws.Range("O11").Value = [double]"3.708144902248023e-14"
$ws.Range("U11").Value = "[44.23652777842656, 56.95498595261816]"
$ws.Range("M12").Value = "[52.15418165874416, 77.28955437378619]"
$ws.Range("N12").Value = [double]"1.63424829224823e-13"
$ws.Range("O12").Value = [double]"1.63424829224823e-13"
$ws.Range("Q12").Value = "[-2.3271056693257726, -1.924579283280233]"
$ws.Range("U12").Value = "[42.86682927076785, 56.0734877589162]"
$ws.Range("Y12").Value = [double]"6.925585585585614"
$ws.Range("Z12").Value = [double]"8.374074074074111"
